$wb = $excel.ActiveWorkbook

# Rename sheets (order corresponds to current sheet order: 1..5)
$wb.Worksheets.Item(1).Name = "GNG_TO-16512555079685133"
$wb.Worksheets.Item(2).Name = "NB_TO-16512555102144668"
$wb.Worksheets.Item(3).Name = "RS_TO-16512555102214544"
$wb.Worksheets.Item(4).Name = "TOL_TO-16512555102794552"
$wb.Worksheets.Item(5).Name = "vSAT_TO-1651255510342451"

# Sheet 1 (GNG_TO)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16512555079325137.csv"
$ws1.Range("B3").Value = "GNG_stims-16512555079515114.csv"
$ws1.Range("B4").Value = "go_stims-1651255507952514.csv"
$ws1.Range("B5").Value = "GNG_stims-16512555079665117.csv"

# Sheet 2 (NB_TO)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16512555086115112.csv"
$ws2.Range("B3").Value = "TB-1651255510063456.csv"
$ws2.Range("B4").Value = "OB-16512555093225105.csv"
$ws2.Range("B5").Value = "ZB-match_0-16512555085375116.csv"
$ws2.Range("B6").Value = "ZB-match_6-16512555083385108.csv"
$ws2.Range("B7").Value = "TB-16512555101934614.csv"
$ws2.Range("B8").Value = "TB-16512555095925097.csv"
$ws2.Range("B9").Value = "ZB-match_2-1651255508196514.csv"
$ws2.Range("B10").Value = "OB-16512555090575132.csv"

# Sheet 3 (RS_TO)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# Sheet 4 (TOL_TO)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16512555102464523.csv"
$ws4.Range("B3").Value = "ZM_stims-1651255510223458.csv"
$ws4.Range("B4").Value = "MM_stims-16512555102624516.csv"
$ws4.Range("B5").Value = "ZM_stims-16512555102464523.csv"
$ws4.Range("B6").Value = "MM_stims-16512555102784538.csv"
$ws4.Range("B7").Value = "ZM_stims-16512555102634528.csv"

# Sheet 5 (vSAT_TO)
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-1651255510283453.csv"
$ws5.Range("B3").Value = "vSAT_stims-16512555103264523.csv"
$ws5.Range("B4").Value = "SAT_stims-16512555102944546.csv"
$ws5.Range("B5").Value = "vSAT_stims-16512555103104599.csv"
